# Update plan import test data:
#  - rename the SHOP QHP sheet to reflect the new "SHOP Q1" naming
#  - flip the "Standard Plan?" values from "Yes" to the new "Y" flag
#  - leave the resulting worksheet selection/active-tab as a SHOP Q1 user
#    would have left it (on D4, with that tab active)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2018_QHP")
$ws1.Name = "SHOP Q1"

$ws1.Range("D2").Value = "Y"
$ws1.Range("D3").Value = "Y"

# Make the renamed sheet the active/selected tab (moves tabSelected from
# 2018_QDP to this sheet) and park the selection on D4.
$ws1.Activate()
$ws1.Range("D4").Select() | Out-Null
